# Add a new "2021" column (column R) to the sheet, mirroring the existing
# 2020 column (Q) formatting, and fill it in with data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell R4, formatted like Q4 (year header)
$ws.Range("Q4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("R4").Value = 2021

# R5 - copy formatting from Q5, then set number format/font/alignment like the source
$ws.Range("Q5").Copy() | Out-Null
$ws.Range("R5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("R5").NumberFormat = "0.00"
$ws.Range("R5").Font.Name = "Times New Roman"
$ws.Range("R5").Font.Bold = $true
$ws.Range("R5").Font.Size = 9
$ws.Range("R5").HorizontalAlignment = -4152   # xlRight
$ws.Range("R5").VerticalAlignment = -4160     # xlTop
$ws.Range("R5").WrapText = $true
$ws.Range("R5").Value = 0.080841202038693286

# R6:R13 - copy formatting from Q6, matching style used by Q6..Q13
$dataRows = 6..13
$values = @{
    6  = 0
    7  = 0
    8  = 0
    9  = 0.2462269049859406
    10 = 0
    11 = 0
    12 = 0
    13 = 0
}
foreach ($r in $dataRows) {
    $ws.Range("Q$r").Copy() | Out-Null
    $ws.Range("R$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("R$r").NumberFormat = "0.00"
    $ws.Range("R$r").Font.Name = "Times New Roman"
    $ws.Range("R$r").Font.Bold = $false
    $ws.Range("R$r").Font.Size = 9
    $ws.Range("R$r").HorizontalAlignment = -4152   # xlRight
    $ws.Range("R$r").VerticalAlignment = -4160     # xlTop
    $ws.Range("R$r").WrapText = $true
    $ws.Range("R$r").Value = $values[$r]
}

# R14 - exactly matches Q14's style (border + number format), just copy format directly
$ws.Range("Q14").Copy() | Out-Null
$ws.Range("R14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("R14").Value = 1.4177257229737372

$excel.CutCopyMode = $false

$ws.Range("T8").Select() | Out-Null
